$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginDataAfterReg")

$ws.Range("A2").Value = "NainaArora6462@testxp.com"
$ws.Range("B2").Value = "Arora@123"
